$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new data row: A7=6, B7=3 (all data collection complete)
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = 3

# Move the active selection down, matching the post-edit cursor position
$ws.Range("B8").Select()
